# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-level holdings, same layout as the
# "2021-Q3" / "2021-Q4" sheets) positioned right before the "总计" (totals)
# sheet, and updates the "总计" sheet with a new leading row summarising the
# 2022-Q1 quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, placed immediately before "总计".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row (bold / centered / bordered, matching the other quarter sheets).
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$headerRange = $q1.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Row index column (A) uses the same style as the header.
$q1.Range("A2").Value = 0
$q1.Range("A2").Font.Bold = $true
$q1.Range("A2").HorizontalAlignment = -4108
$q1.Range("A2").VerticalAlignment = -4160
$q1.Range("A2").Borders.LineStyle = 1

$q1.Range("A3").Value = 1
$q1.Range("A3").Font.Bold = $true
$q1.Range("A3").HorizontalAlignment = -4108
$q1.Range("A3").VerticalAlignment = -4160
$q1.Range("A3").Borders.LineStyle = 1

# Row 2: 003805 / 华安新恒利灵活配置混合A
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "003805"
$q1.Range("C2").Value = "华安新恒利灵活配置混合A"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "5.58"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "24.58"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "0.37"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0206"
$q1.Range("H2").Value = 10

# Row 3: 003806 / 华安新恒利灵活配置混合C
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "003806"
$q1.Range("C3").Value = "华安新恒利灵活配置混合C"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "1.20"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "24.58"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "0.37"
$q1.Range("G3").NumberFormat = "@"
$q1.Range("G3").Value = "0.0044"
$q1.Range("H3").Value = 10

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new leading data row for 2022-Q1,
#    pushing the existing 2021-Q4 / 2021-Q3 rows down by one. The final
#    rows are written with their literal target contents (rather than
#    shifting old values down) to avoid relying on `.Value` reads.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Row 2 (new): 2022-Q1
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.02

# Row 3 (was row 2): 2021-Q4
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.02

# Row 4 (was row 3): 2021-Q3 -- A4 is a brand-new cell, so it needs the same
# bold/centered/bordered format the other index cells in column A carry.
$total.Range("A4").Value = 2
$total.Range("A4").Font.Bold = $true
$total.Range("A4").HorizontalAlignment = -4108
$total.Range("A4").VerticalAlignment = -4160
$total.Range("A4").Borders.LineStyle = 1
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 3
$total.Range("D4").Value = 0.07000000000000001
